# Applies the commit "Add files via upload":
# Renames the folder-structure heading text
#   Rectangular_Outdoor_alternate_UE_arrangement/
# to
#   archive_rectangular_outdoor_alternate_ue_arrangement/
#
# (All other changes in the source diff are Word's automatic removal of
#  proofing-error bookmarks (w:proofErr) and the merging of adjacent runs
#  that share identical formatting -- side effects of re-saving the file
#  in Word, not actual content changes -- so no further edits are needed.)

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Text = "Rectangular_Outdoor_alternate_UE_arrangement"
$find.Replacement.Text = "archive_rectangular_outdoor_alternate_ue_arrangement"
$find.Forward = $true
$find.Wrap = 1            # wdFindContinue
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false

$find.Execute(
    $find.Text,        # FindText
    $true,             # MatchCase
    $false,            # MatchWholeWord
    $false,            # MatchWildcards
    $false,            # MatchSoundsLike
    $false,            # MatchAllWordForms
    $true,             # Forward
    1,                 # Wrap (wdFindContinue)
    $false,            # Format
    $find.Replacement.Text, # ReplaceWith
    2                  # Replace (wdReplaceAll)
) | Out-Null

$d.Save()
